$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-TextValue($addr, $val) {
    $c = $ws.Range($addr)
    $c.NumberFormat = "@"
    $c.Value = $val
    $c.NumberFormat = "General"
    $c.Style = "Normal"
}

Set-TextValue 'D2' '67.023.00'
Set-TextValue 'E2' '  -2.30%  '
Set-TextValue 'D3' '3.751.05'
Set-TextValue 'E3' '  -0.81%  '
Set-TextValue 'D4' '0.999'
Set-TextValue 'E4' '  +0.03%  '
Set-TextValue 'D5' '591.83'
Set-TextValue 'E5' '  -0.85%  '
Set-TextValue 'D6' '165.06'
Set-TextValue 'E6' '  -2.78%  '
Set-TextValue 'D7' '3.749.90'
Set-TextValue 'E7' '  -0.83%  '
Set-TextValue 'E8' '  +0.04%  '
Set-TextValue 'E9' '  -2.02%  '
Set-TextValue 'D10' '0.157'
Set-TextValue 'E10' '  -3.59%  '
Set-TextValue 'D11' '6.36'
Set-TextValue 'E11' '  -2.66%  '
Set-TextValue 'E12' '  -0.88%  '
Set-TextValue 'D13' '0.0000252'
Set-TextValue 'E13' '  -4.93%  '
Set-TextValue 'D14' '35.63'
Set-TextValue 'E14' '  -3.38%  '
Set-TextValue 'D15' '4.378.28'
Set-TextValue 'E15' '  -0.88%  '
Set-TextValue 'D16' '3.760.61'
Set-TextValue 'E16' '  -0.21%  '
Set-TextValue 'D17' '66.987.46'
Set-TextValue 'E17' '  -2.28%  '
Set-TextValue 'D18' '17.65'
Set-TextValue 'E18' '  -2.95%  '
Set-TextValue 'E19' '  -0.10%  '
Set-TextValue 'D20' '6.90'
Set-TextValue 'E20' '  -2.18%  '
Set-TextValue 'D21' '10.44'
Set-TextValue 'E21' '  -4.85%  '
Set-TextValue 'D22' '454.64'
Set-TextValue 'E22' '  -3.09%  '
Set-TextValue 'D23' '0.692'
Set-TextValue 'E23' '  -1.93%  '
Set-TextValue 'D24' '0.0000147'
Set-TextValue 'E24' '  +1.59%  '
Set-TextValue 'D25' '82.69'
Set-TextValue 'E25' '  -2.47%  '
$ws.Range('B26').Value = 'InternetComputer(DFINITY)'
$ws.Range('C26').Value = 'https://coinranking.com/coin/aMNLwaUbY+internetcomputerdfinity-icp'
Set-TextValue 'D26' '11.75'
Set-TextValue 'E26' '  -4.02%  '
$ws.Range('B27').Value = 'Fetch.AI'
$ws.Range('C27').Value = 'https://coinranking.com/coin/AWma-WzFHmKVQ+fetchai-fet'
Set-TextValue 'D27' '2.10'
Set-TextValue 'E27' '  -6.36%  '
Set-TextValue 'E28' '  +0.02%  '
Set-TextValue 'D29' '9.90'
Set-TextValue 'E30' '  -2.05%  '
Set-TextValue 'E31' '  -3.65%  '
Set-TextValue 'D32' '29.50'
Set-TextValue 'E32' '  -2.16%  '
Set-TextValue 'E33' '  -3.27%  '
Set-TextValue 'D34' '9.13'
Set-TextValue 'E34' '  -2.60%  '
Set-TextValue 'D36' '3.704.09'
Set-TextValue 'E36' '  -0.86%  '
Set-TextValue 'D37' '0.0988'
Set-TextValue 'E37' '  -3.22%  '
Set-TextValue 'E38' '  -2.00%  '
Set-TextValue 'D39' '3.25'
Set-TextValue 'E39' '  -7.11%  '
Set-TextValue 'D40' '0.987'
Set-TextValue 'E40' '  -1.80%  '
Set-TextValue 'D41' '5.71'
Set-TextValue 'E41' '  -2.38%  '
Set-TextValue 'E42' '  +0.07%  '
Set-TextValue 'E43' '  -0.02%  '
Set-TextValue 'D44' '43.48'
Set-TextValue 'E44' '  +1.62%  '
Set-TextValue 'D45' '0.296'
Set-TextValue 'E45' '  -4.81%  '
Set-TextValue 'D46' '46.65'
Set-TextValue 'E46' '  +1.84%  '
Set-TextValue 'D47' '8.33'
Set-TextValue 'E47' '  -3.25%  '
Set-TextValue 'D48' '146.84'
Set-TextValue 'E48' '  +1.06%  '
$ws.Range('B49').Value = 'Stacks'
$ws.Range('C49').Value = 'https://coinranking.com/coin/mMPrMcB7+stacks-stx'
Set-TextValue 'D49' '1.81'
Set-TextValue 'E49' '  -8.01%  '
$ws.Range('B50').Value = 'Bittensor'
$ws.Range('C50').Value = 'https://coinranking.com/coin/pgv7xSFi6+bittensor-tao'
Set-TextValue 'D50' '385.57'
Set-TextValue 'E50' '  -3.97%  '
Set-TextValue 'D51' '2.738.05'
Set-TextValue 'E51' '  +1.63%  '
